# Slide 2 formatting fix:
#  - widen the five heading boxes ("Text 3", "Text 6", "Text 9", "Text 12",
#    "Text 15") from 3657600 EMU (288 pt) to 7498080 EMU (590.4 pt) wide
#  - remove the misplaced accent-line shape ("Shape 1")

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Remove the misplaced accent line shape.
$s.Shapes.Item("Shape 1").Delete()

# Widen the heading boxes (3657600 EMU -> 7498080 EMU == 288pt -> 590.4pt).
$headingNames = @("Text 3", "Text 6", "Text 9", "Text 12", "Text 15")
foreach ($name in $headingNames) {
    $shape = $s.Shapes.Item($name)
    $shape.Width = 590.4
}
